# Read sub menu item from excel
# - Submenu sheet gains a new row for a "建物１" (Building 1) object, colored Purple
# - The existing "道" (road) object's color is changed from "blue" to "LightPink"
# - The Submenu sheet becomes the active/selected sheet in the workbook

$wb = $excel.ActiveWorkbook

$wsSubmenu = $wb.Worksheets.Item("Submenu")

# Add the new submenu row (object name / color / flag) before touching the
# existing "blue" cell, so new shared strings are appended in the same
# order the original author typed them in.
$wsSubmenu.Range("A3").Value = "建物１"
$wsSubmenu.Range("B3").Value = "Purple"
$wsSubmenu.Range("C3").Value = 2

# Recolor the existing road entry.
$wsSubmenu.Range("B2").Value = "LightPink"

# Make Submenu the active sheet/tab with B7 selected, matching the
# author's final view state.
$wsSubmenu.Activate() | Out-Null
$wsSubmenu.Range("B7").Select() | Out-Null
